# "updated with latest D models"
#
# Kosten CHF (Tesla Model S electricity price / consumption-factor update):
#   C2 (energy price per 100km baseline figure used in chart "Tesla Model S" column): 72400 -> 73900
#   D4 (consumption factor):                                                           1.8   -> 1.6
#   (B5, B8, C8 all recompute automatically via their existing formulas)
#
# Kosten EUR (same update, EUR sheet):
#   C2: 65300 -> 67900
#   D4: 1.4   -> 1.3
#   (B4, B7, C7 recompute automatically via their existing formulas)

$wb = $excel.ActiveWorkbook

$wsChf = $wb.Worksheets.Item("Kosten CHF")
$wsChf.Range("C2").Value = 73900
$wsChf.Range("D4").Value = 1.6

$wsEur = $wb.Worksheets.Item("Kosten EUR")
$wsEur.Range("C2").Value = 67900
$wsEur.Range("D4").Value = 1.3

# Restore the selection on both sheets to D5 (matches the saved cursor
# position recorded in the workbook after this edit), leaving "Kosten EUR"
# as the active/selected tab just like in the source workbook.
$wsChf.Activate() | Out-Null
$wsChf.Range("D5").Select() | Out-Null

$wsEur.Activate() | Out-Null
$wsEur.Range("D5").Select() | Out-Null
